$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("AppControl")
$ws2 = $wb.Worksheets.Item("smoke")

# --- AppControl (sheet1): add an (empty) hyperlink-styled, wrap-text cell at B25 ---
$ws1.Range("B25").Style = "Hyperlink"
$ws1.Range("B25").WrapText = $true

# --- smoke (sheet2): header A1 renamed "Script Reference" -> "Script_Reference" ---
$ws2.Range("A1").Value = "Script_Reference"

# --- smoke (sheet2): flip the "Email Output" flag from Y to N for rows 18-25 ---
$ws2.Range("B18:B25").Value = "N"

# --- smoke (sheet2): add two new empty rows (26 & 27) with the same border style as B25 ---
$ws2.Range("B25").Copy() | Out-Null
$ws2.Range("B26:B27").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- selection / active-sheet bookkeeping ---
# Make "smoke" the active (selected) sheet/tab, matching the new activeTab/tabSelected state
[void]$ws2.Select()
[void]$ws2.Range("A26:C26").Select()
